# Physics engine position doesn't work 100% so there is an option for the
# client to choose how their position is calculated.
#
# - Mark the "UsePhysics bool" network TODO item (row 26) as DONE, since the
#   feature (letting the client opt in/out of physics-based position calc)
#   has been implemented.
# - Add two new follow-up TODO rows describing the remaining work around
#   client/server start & end position handling now that this option exists.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: "UsePhysics bool" is now finished -> flip status from TODO to DONE
$ws.Range("M26").Value = "DONE"
$ws.Range("M26").Interior.Color = 5287936   # RGB(00,B0,50) green - same as other DONE cells

# Row 27 (new): follow-up task
$ws.Range("K27").Value = "Move start position with avatar"
$ws.Range("M27").Value = "TODO"
$ws.Range("M27").Interior.Color = 255       # RGB(FF,00,00) red - same as other TODO cells

# Row 28 (new): follow-up task
$ws.Range("K28").Value = "Change end position"
$ws.Range("L28").Value = "Without reposistioning avatar"
$ws.Range("M28").Value = "TODO"
$ws.Range("M28").Interior.Color = 255       # RGB(FF,00,00) red - same as other TODO cells

# Match the author's final selection (they had just finished typing M28)
$ws.Range("M28").Select()
